# Update Lgals1-Ptprc TPM-derived NATMI metrics (rows 2-26) with re-computed values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 41.50224066666667
$ws.Range("H2").Value = 124.506722
$ws.Range("I2").Value = 0.05221750323662051
$ws.Range("J2").Value = 0.05952024807977383
$ws.Range("M2").Value = 1.707522
$ws.Range("N2").Value = 5.122566
$ws.Range("O2").Value = 0.002012043481081613
$ws.Range("P2").Value = 0.002016814216842583
$ws.Range("Q2").Value = 70.865988987628
$ws.Range("R2").Value = 637.7939008886519
$ws.Range("S2").Value = 0.0001050638869856003
$ws.Range("T2").Value = 0.0001200412825172853

# Row 3
$ws.Range("G3").Value = 41.50224066666667
$ws.Range("H3").Value = 124.506722
$ws.Range("I3").Value = 0.05221750323662051
$ws.Range("J3").Value = 0.05952024807977383
$ws.Range("O3").Value = 0.000420095518708099
$ws.Range("P3").Value = 0.000421091602904573
$ws.Range("Q3").Value = 14.79614366311622
$ws.Range("R3").Value = 133.165292968046
$ws.Range("S3").Value = 0.00002193633910782993
$ws.Range("T3").Value = 0.0000250634766691898

# Row 4
$ws.Range("G4").Value = 41.50224066666667
$ws.Range("H4").Value = 124.506722
$ws.Range("I4").Value = 0.05221750323662051
$ws.Range("J4").Value = 0.05952024807977383
$ws.Range("M4").Value = 487.9781593333334
$ws.Range("N4").Value = 1463.934478
$ws.Range("O4").Value = 0.5750047580041945
$ws.Range("P4").Value = 0.5763681459167976
$ws.Range("Q4").Value = 20252.18700872902
$ws.Range("R4").Value = 182269.6830785611
$ws.Range("S4").Value = 0.03002531281215622
$ws.Range("T4").Value = 0.03430557503024708

# Row 5
$ws.Range("G5").Value = 41.50224066666667
$ws.Range("H5").Value = 124.506722
$ws.Range("I5").Value = 0.05221750323662051
$ws.Range("J5").Value = 0.05952024807977383
$ws.Range("M5").Value = 6.022401
$ws.Range("N5").Value = 12.044802
$ws.Range("O5").Value = 0.007096443075116684
$ws.Range("P5").Value = 0.004742179585905576
$ws.Range("Q5").Value = 249.943135693174
$ws.Range("R5").Value = 1499.658814159044
$ws.Range("S5").Value = 0.0003705585392433986
$ws.Range("T5").Value = 0.000282255705391939

# Row 6
$ws.Range("G6").Value = 41.50224066666667
$ws.Range("H6").Value = 124.506722
$ws.Range("I6").Value = 0.05221750323662051
$ws.Range("J6").Value = 0.05952024807977383
$ws.Range("M6").Value = 352.5860493333333
$ws.Range("N6").Value = 1057.758148
$ws.Range("O6").Value = 0.415466659920899
$ws.Range("P6").Value = 0.4164517686775497
$ws.Range("Q6").Value = 14633.1110751412
$ws.Range("R6").Value = 131697.9996762708
$ws.Range("S6").Value = 0.02169463165912746
$ws.Range("T6").Value = 0.02478731258494834

# Row 7
$ws.Range("I7").Value = 0.2605968298429562
$ws.Range("J7").Value = 0.2970419303804923
$ws.Range("M7").Value = 1.707522
$ws.Range("N7").Value = 5.122566
$ws.Range("O7").Value = 0.002012043481081613
$ws.Range("P7").Value = 0.002016814216842583
$ws.Range("Q7").Value = 353.664019326576
$ws.Range("R7").Value = 3182.976173939184
$ws.Range("S7").Value = 0.0005243321526760545
$ws.Range("T7").Value = 0.0005990783881897418

# Row 8
$ws.Range("I8").Value = 0.2605968298429562
$ws.Range("J8").Value = 0.2970419303804923
$ws.Range("O8").Value = 0.000420095518708099
$ws.Range("P8").Value = 0.000421091602904573
$ws.Range("S8").Value = 0.0001094755604065629
$ws.Range("T8").Value = 0.0001250818625937901

# Row 9
$ws.Range("I9").Value = 0.2605968298429562
$ws.Range("J9").Value = 0.2970419303804923
$ws.Range("M9").Value = 487.9781593333334
$ws.Range("N9").Value = 1463.934478
$ws.Range("O9").Value = 0.5750047580041945
$ws.Range("P9").Value = 0.5763681459167976
$ws.Range("Q9").Value = 101070.6258387365
$ws.Range("R9").Value = 909635.6325486285
$ws.Range("S9").Value = 0.1498444170805093
$ws.Range("T9").Value = 0.1712055066729508

# Row 10
$ws.Range("I10").Value = 0.2605968298429562
$ws.Range("J10").Value = 0.2970419303804923
$ws.Range("M10").Value = 6.022401
$ws.Range("N10").Value = 12.044802
$ws.Range("O10").Value = 0.007096443075116684
$ws.Range("P10").Value = 0.004742179585905576
$ws.Range("Q10").Value = 1247.366970180408
$ws.Range("R10").Value = 7484.201821082449
$ws.Range("S10").Value = 0.001849310568536407
$ws.Range("T10").Value = 0.001408626178408356

# Row 11
$ws.Range("I11").Value = 0.2605968298429562
$ws.Range("J11").Value = 0.2970419303804923
$ws.Range("M11").Value = 352.5860493333333
$ws.Range("N11").Value = 1057.758148
$ws.Range("O11").Value = 0.415466659920899
$ws.Range("P11").Value = 0.4164517686775497
$ws.Range("Q11").Value = 73028.0484618676
$ws.Range("R11").Value = 657252.4361568084
$ws.Range("S11").Value = 0.1082692944808279
$ws.Range("T11").Value = 0.1237036372783496

# Row 12
$ws.Range("G12").Value = 104.6648203333333
$ws.Range("H12").Value = 313.994461
$ws.Range("I12").Value = 0.1316877235234609
$ws.Range("J12").Value = 0.1501045719796146
$ws.Range("M12").Value = 1.707522
$ws.Range("N12").Value = 5.122566
$ws.Range("O12").Value = 0.002012043481081613
$ws.Range("P12").Value = 0.002016814216842583
$ws.Range("Q12").Value = 178.717483345214
$ws.Range("R12").Value = 1608.457350106926
$ws.Range("S12").Value = 0.0002649614256538574
$ws.Range("T12").Value = 0.0003027330347815577

# Row 13
$ws.Range("G13").Value = 104.6648203333333
$ws.Range("H13").Value = 313.994461
$ws.Range("I13").Value = 0.1316877235234609
$ws.Range("J13").Value = 0.1501045719796146
$ws.Range("O13").Value = 0.000420095518708099
$ws.Range("P13").Value = 0.000421091602904573
$ws.Range("Q13").Value = 37.31450864459144
$ws.Range("R13").Value = 335.830577801323
$ws.Range("S13").Value = 0.00005532142252107705
$ws.Range("T13").Value = 0.00006320777481820078

# Row 14
$ws.Range("G14").Value = 104.6648203333333
$ws.Range("H14").Value = 313.994461
$ws.Range("I14").Value = 0.1316877235234609
$ws.Range("J14").Value = 0.1501045719796146
$ws.Range("M14").Value = 487.9781593333334
$ws.Range("N14").Value = 1463.934478
$ws.Range("O14").Value = 0.5750047580041945
$ws.Range("P14").Value = 0.5763681459167976
$ws.Range("Q14").Value = 51074.14637321405
$ws.Range("R14").Value = 459667.3173589264
$ws.Range("S14").Value = 0.07572106759673093
$ws.Range("T14").Value = 0.08651549384552498

# Row 15
$ws.Range("G15").Value = 104.6648203333333
$ws.Range("H15").Value = 313.994461
$ws.Range("I15").Value = 0.1316877235234609
$ws.Range("J15").Value = 0.1501045719796146
$ws.Range("M15").Value = 6.022401
$ws.Range("N15").Value = 12.044802
$ws.Range("O15").Value = 0.007096443075116684
$ws.Range("P15").Value = 0.004742179585905576
$ws.Range("Q15").Value = 630.333518640287
$ws.Range("R15").Value = 3782.001111841722
$ws.Range("S15").Value = 0.0009345144336759449
$ws.Range("T15").Value = 0.0007118228369928226

# Row 16
$ws.Range("G16").Value = 104.6648203333333
$ws.Range("H16").Value = 313.994461
$ws.Range("I16").Value = 0.1316877235234609
$ws.Range("J16").Value = 0.1501045719796146
$ws.Range("M16").Value = 352.5860493333333
$ws.Range("N16").Value = 1057.758148
$ws.Range("O16").Value = 0.415466659920899
$ws.Range("P16").Value = 0.4164517686775497
$ws.Range("Q16").Value = 36903.35550551313
$ws.Range("R16").Value = 332130.1995496182
$ws.Range("S16").Value = 0.05471185864487912
$ws.Range("T16").Value = 0.06251131448749708

# Row 17
$ws.Range("G17").Value = 292.548645
$ws.Range("H17").Value = 585.0972899999999
$ws.Range("I17").Value = 0.3680803631748439
$ws.Range("J17").Value = 0.2797048648634679
$ws.Range("M17").Value = 1.707522
$ws.Range("N17").Value = 5.122566
$ws.Range("O17").Value = 0.002012043481081613
$ws.Range("P17").Value = 0.002016814216842583
$ws.Range("Q17").Value = 499.5332474076899
$ws.Range("R17").Value = 2997.19948444614
$ws.Range("S17").Value = 0.0007405936952400973
$ws.Range("T17").Value = 0.0005641127479766757

# Row 18
$ws.Range("G18").Value = 292.548645
$ws.Range("H18").Value = 585.0972899999999
$ws.Range("I18").Value = 0.3680803631748439
$ws.Range("J18").Value = 0.2797048648634679
$ws.Range("O18").Value = 0.000420095518708099
$ws.Range("P18").Value = 0.000421091602904573
$ws.Range("Q18").Value = 104.297785139745
$ws.Range("R18").Value = 625.7867108384698
$ws.Range("S18").Value = 0.0001546289110942015
$ws.Range("T18").Value = 0.0001177813698855647

# Row 19
$ws.Range("G19").Value = 292.548645
$ws.Range("H19").Value = 585.0972899999999
$ws.Range("I19").Value = 0.3680803631748439
$ws.Range("J19").Value = 0.2797048648634679
$ws.Range("M19").Value = 487.9781593333334
$ws.Range("N19").Value = 1463.934478
$ws.Range("O19").Value = 0.5750047580041945
$ws.Range("P19").Value = 0.5763681459167976
$ws.Range("Q19").Value = 142757.3493025608
$ws.Range("R19").Value = 856544.0958153646
$ws.Range("S19").Value = 0.2116479601534471
$ws.Range("T19").Value = 0.1612129743652654

# Row 20
$ws.Range("G20").Value = 292.548645
$ws.Range("H20").Value = 585.0972899999999
$ws.Range("I20").Value = 0.3680803631748439
$ws.Range("J20").Value = 0.2797048648634679
$ws.Range("M20").Value = 6.022401
$ws.Range("N20").Value = 12.044802
$ws.Range("O20").Value = 0.007096443075116684
$ws.Range("P20").Value = 0.004742179585905576
$ws.Range("Q20").Value = 1761.845252196645
$ws.Range("R20").Value = 7047.38100878658
$ws.Range("S20").Value = 0.002612061344338555
$ws.Range("T20").Value = 0.001326410700234015

# Row 21
$ws.Range("G21").Value = 292.548645
$ws.Range("H21").Value = 585.0972899999999
$ws.Range("I21").Value = 0.3680803631748439
$ws.Range("J21").Value = 0.2797048648634679
$ws.Range("M21").Value = 352.5860493333333
$ws.Range("N21").Value = 1057.758148
$ws.Range("O21").Value = 0.415466659920899
$ws.Range("P21").Value = 0.4164517686775497
$ws.Range("Q21").Value = 103148.5709783698
$ws.Range("R21").Value = 618891.4258702188
$ws.Range("S21").Value = 0.1529251190707239
$ws.Range("T21").Value = 0.1164835856801062

# Row 22
$ws.Range("G22").Value = 148.958664
$ws.Range("H22").Value = 446.875992
$ws.Range("I22").Value = 0.1874175802221185
$ws.Range("J22").Value = 0.2136283846966514
$ws.Range("M22").Value = 1.707522
$ws.Range("N22").Value = 5.122566
$ws.Range("O22").Value = 0.002012043481081613
$ws.Range("P22").Value = 0.002016814216842583
$ws.Range("Q22").Value = 254.350195870608
$ws.Range("R22").Value = 2289.151762835472
$ws.Range("S22").Value = 0.0003770923205260037
$ws.Range("T22").Value = 0.000430848763377323

# Row 23
$ws.Range("G23").Value = 148.958664
$ws.Range("H23").Value = 446.875992
$ws.Range("I23").Value = 0.1874175802221185
$ws.Range("J23").Value = 0.2136283846966514
$ws.Range("O23").Value = 0.000420095518708099
$ws.Range("P23").Value = 0.000421091602904573
$ws.Range("Q23").Value = 53.105898790184
$ws.Range("R23").Value = 477.953089111656
$ws.Range("S23").Value = 0.00007873328557842761
$ws.Range("T23").Value = 0.00008995711893782768

# Row 24
$ws.Range("G24").Value = 148.958664
$ws.Range("H24").Value = 446.875992
$ws.Range("I24").Value = 0.1874175802221185
$ws.Range("J24").Value = 0.2136283846966514
$ws.Range("M24").Value = 487.9781593333334
$ws.Range("N24").Value = 1463.934478
$ws.Range("O24").Value = 0.5750047580041945
$ws.Range("P24").Value = 0.5763681459167976
$ws.Range("Q24").Value = 72688.57467547247
$ws.Range("R24").Value = 654197.1720792522
$ws.Range("S24").Value = 0.1077660003613509
$ws.Range("T24").Value = 0.1231285960028093

# Row 25
$ws.Range("G25").Value = 148.958664
$ws.Range("H25").Value = 446.875992
$ws.Range("I25").Value = 0.1874175802221185
$ws.Range("J25").Value = 0.2136283846966514
$ws.Range("M25").Value = 6.022401
$ws.Range("N25").Value = 12.044802
$ws.Range("O25").Value = 0.007096443075116684
$ws.Range("P25").Value = 0.004742179585905576
$ws.Range("Q25").Value = 897.0888070322641
$ws.Range("R25").Value = 5382.532842193585
$ws.Range("S25").Value = 0.001329998189322378
$ws.Range("T25").Value = 0.001013064164878443

# Row 26
$ws.Range("G26").Value = 148.958664
$ws.Range("H26").Value = 446.875992
$ws.Range("I26").Value = 0.1874175802221185
$ws.Range("J26").Value = 0.2136283846966514
$ws.Range("M26").Value = 352.5860493333333
$ws.Range("N26").Value = 1057.758148
$ws.Range("O26").Value = 0.415466659920899
$ws.Range("P26").Value = 0.4164517686775497
$ws.Range("Q26").Value = 52520.74685373141
$ws.Range("R26").Value = 472686.7216835828
$ws.Range("S26").Value = 0.0778657560653407
$ws.Range("T26").Value = 0.08896591864664845

Write-Host "Updated 278 cells across rows 2-26."
